$d = $word.ActiveDocument

$pairs = @(
    @("426×5=", "660×5="),
    @("541×4=", "784×3="),
    @("575×7=", "667×8="),
    @("832×9=", "269×5="),
    @("924×2=", "197×4="),
    @("211×7=", "282×3="),
    @("843×6=", "212×2="),
    @("992×2=", "868×6="),
    @("880×8=", "201×4="),
    @("456×9=", "254×4="),
    @("631×7=", "818×8="),
    @("820×4=", "387×3="),
    @("239×6=", "603×4="),
    @("502×6=", "499×6="),
    @("762×9=", "896×2="),
    @("416×2=", "695×3="),
    @("221×6=", "168×9="),
    @("533×9=", "445×5="),
    @("421×4=", "879×6="),
    @("545×8=", "494×6="),
    @("142×4=", "736×8="),
    @("304×7=", "187×4="),
    @("364×5=", "813×4="),
    @("864×6=", "389×2="),
    @("645×7=", "632×4=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}
